# "Added project report and ppt"
# The crime-report backing data (cols R:U on Sheet1) is refreshed: the old
# single "R" column of numbers is replaced by three new columns S/T/U with
# updated figures for several neighborhoods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old column-R figures (rows 13-24) and the two stray old T values
# that don't carry over to the refreshed data.
$ws.Range("R13:R24").ClearContents()
$ws.Range("T19").ClearContents()
$ws.Range("T22").ClearContents()

# Del Mar Heights (row 14)
$ws.Range("S14").Value = 20
$ws.Range("T14").Value = 47
$ws.Range("U14").Value = 9

# Rancho Penasquitos (row 15)
$ws.Range("S15").Value = 97
$ws.Range("T15").Value = 304
$ws.Range("U15").Value = 50

# Bay Park (row 16)
$ws.Range("S16").Value = 22
$ws.Range("T16").Value = 33
$ws.Range("U16").Value = 15

# Normal Heights (row 17)
$ws.Range("S17").Value = 13
$ws.Range("T17").Value = 31
$ws.Range("U17").Value = 7

# Bay Terraces (row 21)
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 14
$ws.Range("U21").Value = 4

# Nudge the saved window position/selection to match the refreshed view.
$win = $excel.ActiveWindow
$win.Left = 15000
$win.Top = 900
$ws.Range("V29").Select() | Out-Null
